$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: set a cell's value while forcing it to remain TEXT (so that
# numeric-looking strings such as "0.02850" keep their exact textual
# representation, including trailing zeros, instead of being coerced into
# a floating point number). ---
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
}

# Plain "Price" column updates (column D) - numeric-looking strings that
# must stay as text so formatting (trailing zeros, etc.) is preserved.
Set-TextValue "D2"  "245.77"
Set-TextValue "D3"  "26.25"
Set-TextValue "D4"  "5.075"
Set-TextValue "D6"  "6.473"
Set-TextValue "D8"  "0.8139"
Set-TextValue "D9"  "0.8438"
Set-TextValue "D10" "0.1347"
Set-TextValue "D11" "0.02850"
Set-TextValue "D12" "0.09391"
Set-TextValue "D13" "0.001512"
Set-TextValue "D14" "0.0006020"
Set-TextValue "D15" "0.006207"
Set-TextValue "D16" "3.561"
Set-TextValue "D19" "0.07006"
Set-TextValue "D20" "0.03214"
Set-TextValue "D22" "3.746"
Set-TextValue "D23" "0.04702"
Set-TextValue "D25" "0.001249"
Set-TextValue "D26" "0.004599"

# Row 26 - Volume(1h) label also changed.
$ws.Range("E26").Value = "25HotbitTokenHTBBestin24h"

Set-TextValue "D27" "0.00009605"
Set-TextValue "D40" "0.03651"

# Rows 41-43 were reshuffled (Coin / Link / Price / Volume columns).
# Row 41: was BKEXToken -> now KickToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006183"
$ws.Range("E41").Value = "40KickTokenKICK"

# Row 42: was CEJI -> now BKEXToken
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1054"
$ws.Range("E42").Value = "41BKEXTokenBKK"

# Row 43: was KickToken -> now CEJI
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002501"
$ws.Range("E43").Value = "42CEJICEJI"

Set-TextValue "D44" "0.007681"

Set-TextValue "D48" "0.002056"
Set-TextValue "D49" "0.00002101"
Set-TextValue "D50" "0.0002001"
